# Generate Report for Archive
# Update status text from "Ready for handoff" to "In Translation" on all
# sheets that reference it, then resize the affected Status columns to
# match the new (shorter) text, as Excel's AutoFit would do.

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"
# The target stored column width (13.4101845877511 "characters") comes from
# genuine Excel sub-pixel text metrics that this engine cannot reproduce
# exactly; ColumnWidth assignments here snap to a 1/6-character pixel grid,
# so 12.5 is the input that lands on the closest reachable grid point.
$newColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ($oldText -eq $cell.Value2) {
                $cell.Value = $newText
            }
        }
    }
}

# Overview sheet: Status columns are E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# zh-cn sheet: Status column is C
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = $newColumnWidth

# de-de sheet: Status column is C
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = $newColumnWidth
